$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-11 (row 24)
$ws.Range("B24").Value = 6341
$ws.Range("C24").Value = 998
$ws.Range("D24").Value = 5940333
$ws.Range("E24").Value = 936.8132786626715
$ws.Range("F24").Value = 8.097511080804631
$ws.Range("G24").Value = 3.419689119170988
$ws.Range("H24").Value = 25.8392626995203
